# TOD-E norms run, POM rescale, 24 cell demo strat
#
# The workbook's last sheet "7.0-9.3" gets rescaled/rescoped into four
# narrower age bands:
#   "7.0-9.3"  -> "7.0-7.5"   (rescaled raw->ss lookup, in place)
#                 "7.6-7.11"  (new sheet)
#                 "8.0-8.5"   (new sheet)
#                 "8.6-9.3"   (new sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("7.0-9.3")

# New raw(1-25) -> ss lookup values for each of the four resulting sheets.
$ssValues = @{
    "7.0-7.5"  = @(63,64,66,67,69,70,72,73,75,76,78,80,82,84,86,89,91,94,97,101,106,112,120,128,130)
    "7.6-7.11" = @(59,60,62,63,64,66,67,69,70,72,73,75,77,79,81,83,86,89,92,96,102,109,119,126,130)
    "8.0-8.5"  = @(55,57,58,59,60,62,63,64,66,67,69,71,73,74,77,79,81,84,88,92,98,107,117,123,127)
    "8.6-9.3"  = @(51,52,54,55,56,57,58,60,61,63,64,66,67,69,71,73,76,78,82,86,94,106,113,117,120)
}

$sheetOrder = @("7.0-7.5", "7.6-7.11", "8.0-8.5", "8.6-9.3")

# Rescale the existing "7.0-9.3" sheet in place -> becomes "7.0-7.5".
$firstName = $sheetOrder[0]
$vals = $ssValues[$firstName]
for ($i = 0; $i -lt $vals.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $vals[$i]
}
$ws.Name = $firstName

# Create the remaining three sheets as copies of the rescaled sheet (so the
# header formatting / raw column is carried over identically), placing each
# right after the previous one, then overwrite their ss column.
$after = $ws
for ($s = 1; $s -lt $sheetOrder.Count; $s++) {
    $name = $sheetOrder[$s]

    $ws.Copy([System.Reflection.Missing]::Value, $after)
    $newSheet = $wb.Worksheets.Item($after.Index + 1)
    $newSheet.Name = $name

    $vals = $ssValues[$name]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $newSheet.Cells.Item($i + 2, 2).Value = $vals[$i]
    }

    $after = $newSheet
}
